# Auto-generated edit script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $val) {
    # Force text interpretation so numeric-looking strings (e.g. "53.00", "0.818")
    # are preserved exactly instead of being coerced into floating point numbers,
    # then restore the default "Normal" style so no stray style index is left behind.
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextCell "D2" "43.540.92"
$ws.Range("E2").Value = "  +1.04%  "
Set-TextCell "D3" "2.375.08"
$ws.Range("E3").Value = "  +3.20%  "
Set-TextCell "D4" "0.999"
$ws.Range("E4").Value = "  -0.16%  "
Set-TextCell "D5" "310.93"
$ws.Range("E5").Value = "  +0.25%  "
Set-TextCell "D6" "105.31"
$ws.Range("E6").Value = "  +4.45%  "
$ws.Range("E7").Value = "  -0.50%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  +3.10%  "
Set-TextCell "D10" "36.48"
$ws.Range("E10").Value = "  +0.80%  "
Set-TextCell "D11" "53.00"
$ws.Range("E11").Value = "  +1.88%  "
Set-TextCell "D12" "0.0817"
$ws.Range("E12").Value = "  -0.33%  "
$ws.Range("E13").Value = "  -0.41%  "
$ws.Range("E14").Value = "  +1.20%  "
Set-TextCell "D15" "2.743.51"
$ws.Range("E15").Value = "  +3.33%  "
Set-TextCell "D16" "15.74"
$ws.Range("E16").Value = "  +6.00%  "
Set-TextCell "D17" "2.371.20"
$ws.Range("E17").Value = "  +2.42%  "
Set-TextCell "D18" "0.818"
$ws.Range("E18").Value = "  +1.64%  "
Set-TextCell "D19" "43.541.37"
$ws.Range("E19").Value = "  +1.07%  "
Set-TextCell "D20" "12.05"
$ws.Range("E20").Value = "  -3.92%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextCell "D21" "6.33"
$ws.Range("E21").Value = "  +4.55%  "
$ws.Range("B22").Value = "ShibaInu"
$ws.Range("C22").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextCell "D22" "0.0₃0931"
$ws.Range("E22").Value = "  +1.48%  "
Set-TextCell "D23" "68.55"
$ws.Range("E23").Value = "  +0.94%  "
Set-TextCell "D24" "243.79"
$ws.Range("E24").Value = "  +1.39%  "
Set-TextCell "D25" "2.07"
$ws.Range("E25").Value = "  +2.42%  "
Set-TextCell "D26" "2.64"
$ws.Range("E26").Value = "  +1.08%  "
$ws.Range("E27").Value = "  +0.17%  "
Set-TextCell "D28" "25.95"
$ws.Range("E28").Value = "  +8.62%  "
$ws.Range("B29").Value = "InjectiveProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextCell "D29" "36.84"
$ws.Range("E29").Value = "  -4.09%  "
$ws.Range("B30").Value = "Cosmos"
$ws.Range("C30").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextCell "D30" "9.64"
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextCell "D31" "2.11"
$ws.Range("E31").Value = "  -2.08%  "
Set-TextCell "D32" "162.54"
$ws.Range("E32").Value = "  -1.37%  "
Set-TextCell "D33" "5.32"
$ws.Range("E33").Value = "  +0.06%  "
Set-TextCell "D34" "0.998"
$ws.Range("E34").Value = "  -0.29%  "
Set-TextCell "D35" "18.37"
$ws.Range("E35").Value = "  +3.28%  "
Set-TextCell "D36" "3.16"
$ws.Range("E36").Value = "  +0.60%  "
$ws.Range("E37").Value = "  +6.47%  "
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextCell "D38" "0.0744"
$ws.Range("E38").Value = "  +0.50%  "
$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextCell "D39" "1.96"
$ws.Range("E39").Value = "  +6.41%  "
Set-TextCell "D40" "4.66"
$ws.Range("E40").Value = "  +11.97%  "
Set-TextCell "D41" "0.106"
$ws.Range("E41").Value = "  +1.15%  "
$ws.Range("E42").Value = "  -0.54%  "
Set-TextCell "D43" "2.42"
$ws.Range("E43").Value = "  +5.10%  "
Set-TextCell "D44" "20.38"
$ws.Range("E44").Value = "  +5.06%  "
Set-TextCell "D45" "2.009.92"
$ws.Range("E45").Value = "  +2.23%  "
Set-TextCell "D46" "0.0293"
$ws.Range("E46").Value = "  +1.20%  "
Set-TextCell "D47" "3.17"
$ws.Range("E47").Value = "  +4.80%  "
Set-TextCell "D48" "10.44"
$ws.Range("E48").Value = "  +5.98%  "
Set-TextCell "D49" "58.34"
$ws.Range("E49").Value = "  +6.44%  "
Set-TextCell "D50" "2.92"
$ws.Range("E50").Value = "  -2.90%  "
$ws.Range("E51").Value = "  +2.29%  "
